$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 28.51638556701447
